$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.531.10"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "1.906.42"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'338.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.13%  "

$ws.Range("D6").Value = "'1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4760"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.21%  "

$ws.Range("D8").Value = "'0.4002"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.80%  "

$ws.Range("D9").Value = "'0.08032"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.50%  "

$ws.Range("D10").Value = "'0.9916"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("D11").Value = "'23.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").Value = "1.894.32"
$ws.Range("E12").Value = "  -1.31%  "

$ws.Range("D13").Value = "'5.922"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.56%  "

$ws.Range("D14").Value = "'7.113"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.83%  "

$ws.Range("D15").Value = "'89.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.57%  "

$ws.Range("D16").Value = "'0.06833"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "'1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("D20").Value = "'1.008"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").Value = "29.550.29"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").Value = "'5.506"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.70%  "

$ws.Range("D23").Value = "'11.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.95%  "

$ws.Range("D24").Value = "'2.150"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").Value = "2.144.09"
$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("D26").Value = "'157.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("D27").Value = "'6.498"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.53%  "

$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("D29").Value = "'2.053"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.90%  "

$ws.Range("D30").Value = "'119.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.08%  "

$ws.Range("D31").Value = "'0.9949"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").Value = "'0.09520"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("D33").Value = "'5.476"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.97%  "

$ws.Range("D34").Value = "'3.543"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").Value = "'1.386"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.97%  "

$ws.Range("D36").Value = "'0.06464"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.78%  "

$ws.Range("D37").Value = "'0.02241"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("E38").Value = "  +1.12%  "

$ws.Range("D39").Value = "'0.5817"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.79%  "

$ws.Range("D40").Value = "'10.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.09%  "

$ws.Range("D41").Value = "'7.749"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.25%  "

$ws.Range("D42").Value = "'0.1821"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("D43").Value = "'2.457"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.09%  "

$ws.Range("D44").Value = "'1.270"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").Value = "'0.07415"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.43%  "

$ws.Range("D46").Value = "'12.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.23%  "

$ws.Range("D47").Value = "'0.5481"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("D48").Value = "'1.947"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.38%  "

$ws.Range("D49").Value = "'116.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.91%  "

$ws.Range("D50").Value = "'2.374"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("D51").Value = "'71.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.68%  "

